$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "201124"

# Update the "Create Date" value in B2 (shared string)
$ws.Range("B2").Value = "2024-02-03T14:56:59.270677"

# Update Quantity (C2) and Total Price (D2)
$ws.Range("C2").Value = 1.0
$ws.Range("D2").Value = 78000.0
